# Power_Storage.xlsx - "Increase MaxInvest Storage Adapt Szenarios Existing Units"
#
# Bumps the MaxInvest (column S) value for the existing storage units from
# 8 to 15 MW across rows 7-11, and bumps the ExisUnits (column E) count for
# row 10 from 31 to 33. Also updates the active selection on the sheet to
# reflect where the user was last working (S8:S11 in the frozen bottom-left
# pane).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power Storage")
$ws.Activate()

# MaxInvest column (S) - existing units across the storage technologies
$ws.Range("S7").Value = 15
$ws.Range("S8").Value = 15
$ws.Range("S9").Value = 15
$ws.Range("S10").Value = 15
$ws.Range("S11").Value = 15

# ExisUnits column (E) for row 10
$ws.Range("E10").Value = 33

# Leave the selection where the edits were made (bottom-left frozen pane)
$ws.Range("S8:S11").Select()
